$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 464
$ws1.Range("F9").Value = 1753
$ws1.Range("F10").Value = 375
$ws1.Range("F15").Value = 12855
$ws1.Range("F16").Value = 12842
$ws1.Range("F22").Value = 580
$ws1.Range("F27").Value = 75
$ws1.Range("F28").Value = 255

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 82

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 464
$ws4.Range("F14").Value = 1753
$ws4.Range("F15").Value = 375
$ws4.Range("F21").Value = 12855
$ws4.Range("F22").Value = 12842
$ws4.Range("F28").Value = 580
$ws4.Range("F37").Value = 75
$ws4.Range("F38").Value = 255
$ws4.Range("F40").Value = 82
